$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: Nitra Nylon Tubing (flexible tubing) ---
# Fill order matters for shared-string insertion order: Name, Description,
# Part Number (D) before Quantity (C), then Cost.
$ws.Range("A17").Value = "Nitra Nylon Tubing"
$ws.Range("B17").Value = "Flexible Tubing"
$ws.Range("D17").Value = "n532BLU100"
$ws.Range("C17").Value = "200 ft"
$ws.Range("E17").Value = 21
$ws.Range("E17").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# --- Row 18: Bimba Air Cylinder ---
$ws.Range("A18").Value = "Bimba Air Cylinder"
$ws.Range("B18").Value = "2 IN Bore Diameter Dual Acting Cylinders with position feedback"
$ws.Range("C18").Value = 8
$ws.Range("E18").Value = 2831.2
$ws.Range("E18").NumberFormat = "$#,##0.00_);[Red](""$""#,##0.00)"

# --- Column width adjustments ---
# Column E (Cost) widened to fit the new currency values.
$ws.Columns("E").ColumnWidth = 10
# Column I picked up a stray width (no data), matching a leftover from the
# editing session that produced this workbook.
$ws.Columns("I").ColumnWidth = 13.3

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection / active cell ---
[void]$ws.Range("I23").Select()
